$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H holds the "SnippetID" values; update them in place (row -> new value)
$ws.Range("H2").Value  = "pRc8"
$ws.Range("H3").Value  = "jYRk"
$ws.Range("H4").Value  = "44Bu"
$ws.Range("H5").Value  = "fS1D"
$ws.Range("H6").Value  = "2qVk"
$ws.Range("H7").Value  = "Mt0v"
$ws.Range("H8").Value  = "Mt0v"
$ws.Range("H9").Value  = "Mt0v"
$ws.Range("H10").Value = "Mt0v"
$ws.Range("H11").Value = "Mt0v"
$ws.Range("H12").Value = "9jkf"
$ws.Range("H13").Value = "7ugC"
$ws.Range("H14").Value = "TDsX"
$ws.Range("H15").Value = "qcyb"
$ws.Range("H16").Value = "c8Ox"
$ws.Range("H17").Value = "c8Ox"
$ws.Range("H18").Value = "J12J"
$ws.Range("H19").Value = "Tt3X"
$ws.Range("H20").Value = "VWFS"
$ws.Range("H21").Value = "VWFS"
$ws.Range("H22").Value = "lAhm"
$ws.Range("H23").Value = "bE3k"
$ws.Range("H24").Value = "sU8t"
$ws.Range("H25").Value = "QCZl"
$ws.Range("H26").Value = "NB3S"
$ws.Range("H27").Value = "v3O6"
$ws.Range("H28").Value = "THo6"
$ws.Range("H29").Value = "4q8u"
